# Compiled_Data_Handwritten.xlsx edit:
#  - Sheet1: columns D/F/H become literal (2-decimal) input values, and
#    columns C/E/G become CEILING(B/D,1) / CEILING(B/F,1) / CEILING(B/H,1)
#    formulas (row 3 entered individually, rows 4:25 entered as one range
#    so Excel records them as a shared-formula group).
#  - Sheet1 column A (rows 3-14) relabeled to the new "chN" file names.
#  - Sheet2 keeps its original formulas/values but they become shared
#    formula groups spanning rows 3:25 (same values, just re-entered as a
#    single range so they can be recorded as shared formulas).
#  - Sheet1 selection moves to G3:H25.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# New "column D/F/H" input values (rows 3..25), 2 decimal places in the sheet.
$D = @(1.357,1.311,1.342,1.654,1.573,1.621,0.304,0.313,0.313,0.756,0.755,0.755,1.076,0.754,0.49,0.887,0.887,3.761,0.918,0.887,0.887,0.905,0.91)
$F = @(2.083,1.944,2.05,3.302,2.893,3.195,1.691,1.704,1.68,0.832,0.833,0.834,9.156,1.923,2.49,1.654,1.336,7.555,4.58,1.818,1.054,0.731,0.695)
$H = @(0.947,0.888,0.934,1.261,1.123,1.221,0.811,1.0,0.807,1.0,1.007,1.005,3.523,0.858,1.0,0.805,0.78,3.984,1.84,0.84,1.0,0.664,0.651)

# New shared-string labels for A3:A14 (data0_ch1.raw .. data3_ch3.raw).
$labels = @("data0_ch1.raw","data0_ch2.raw","data0_ch3.raw","data1_ch1.raw","data1_ch2.raw","data1_ch3.raw","data2_ch1.raw","data2_ch2.raw","data2_ch3.raw","data3_ch1.raw","data3_ch2.raw","data3_ch3.raw")

for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws1.Cells.Item(3 + $i, 1).Value = $labels[$i]
}

for ($i = 0; $i -lt $D.Length; $i++) {
    $r = 3 + $i
    $ws1.Cells.Item($r, 4).Value = $D[$i]
    $ws1.Cells.Item($r, 6).Value = $F[$i]
    $ws1.Cells.Item($r, 8).Value = $H[$i]
}

# C/E/G are now derived from D/F/H via CEILING(..,1). Row 3 is entered on
# its own (matches the diff: C3/E3/G3 are plain, non-shared formulas);
# rows 4:25 are entered as a single range each so they become one shared
# formula group, same as the source workbook.
$ws1.Range("C3").Formula = "=CEILING(B3/D3,1)"
$ws1.Range("E3").Formula = "=CEILING(B3/F3,1)"
$ws1.Range("G3").Formula = "=CEILING(B3/H3,1)"

$ws1.Range("C4:C25").Formula = "=CEILING(B4/D4,1)"
$ws1.Range("E4:E25").Formula = "=CEILING(B4/F4,1)"
$ws1.Range("G4:G25").Formula = "=CEILING(B4/H4,1)"

# D/F/H display with 2 decimal places.
$ws1.Range("D3:D25").NumberFormat = "0.00"
$ws1.Range("F3:F25").NumberFormat = "0.00"
$ws1.Range("H3:H25").NumberFormat = "0.00"

# New selection on Sheet1.
$ws1.Range("G3:H25").Select()

# Sheet2: same formulas/values as before, but re-entered as whole ranges
# (rows 3:25) so they are recorded as shared-formula groups.
$ws2.Range("D3:D25").Formula = "=B3/C3"
$ws2.Range("E3:E25").Formula = "=CEILING(B3/F3,1)"
$ws2.Range("H3:H25").Formula = "=B3/G3"
